# 2021 Excess Mortality Update
# - Revise B17 (2020 Q4) and B18 (2021 Q1) death-rate figures
# - Append three new quarters of data: 2021 Q2, 2021 Q3, 2021 Q4
# - Grow Table3 so the new rows are included in the table/autofilter

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing figures -------------------------------------------------
$ws.Range("B17").Value = 194
$ws.Range("B18").Value = 214.5

# --- Append the new quarters --------------------------------------------------
$ws.Range("A19").Value = "2021 Q2"
$ws.Range("B19").Value = 142.7

$ws.Range("A20").Value = "2021 Q3"
$ws.Range("B20").Value = 160.7

$ws.Range("A21").Value = "2021 Q4"
$ws.Range("B21").Value = 167.4

# --- Grow the table / autofilter so it covers the new rows -------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B21"))
